$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, shifting rows 29-105 down to 30-106
$ws.Rows("29").Insert()

# Populate the newly inserted row 29 with the new record
$ws.Cells.Item(29, 1).Value2 = 10
$ws.Cells.Item(29, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(29, 3).Value2 = "La Araucanía"
$ws.Cells.Item(29, 4).Value2 = 44525
$ws.Cells.Item(29, 5).Value2 = 9
$ws.Cells.Item(29, 6).Value2 = 100112012
$ws.Cells.Item(29, 7).Value2 = "Espinaca"
$ws.Cells.Item(29, 8).Value2 = "Sin especificar"
$ws.Cells.Item(29, 9).Value2 = "Primera"
$ws.Cells.Item(29, 10).Value2 = 50
$ws.Cells.Item(29, 11).Value2 = 8000
$ws.Cells.Item(29, 12).Value2 = 8000
$ws.Cells.Item(29, 13).Value2 = 8000
$ws.Cells.Item(29, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(29, 15).Value2 = "Región de La Araucanía"
$ws.Cells.Item(29, 16).Value2 = 2667
$ws.Cells.Item(29, 17).Value2 = 3
$ws.Cells.Item(29, 18).Value2 = "Hortaliza"
